$d = $word.ActiveDocument

# 1. Update the title text
$d.Content.Find.Execute("Child Care Home Page", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Provider Details Page", 2)

# 2. Update the user story text
$d.Content.Find.Execute( `
    "As a caseworker or parent, I need a State of Mississippi Home page so that I can begin my search for child care providers.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "As a caseworker or parent, I would like to see more details about the provider so that I can make an informed decision.", 2)

# 3. Move the "_GoBack" bookmark from the "FORM FIELDS" paragraph to the
#    empty paragraph that precedes it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $d.Paragraphs.Item(8)
$d.Bookmarks.Add("_GoBack", $target.Range)
